$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each D-column price is stored as text (inlineStr) in the original workbook.
# Setting NumberFormat to "@" (Text) before assigning the new value keeps Excel
# from auto-coercing the numeric-looking string into a real number, and
# ClearFormats() afterwards removes the now-unneeded explicit format so the
# cell keeps its original (default) style.
$priceUpdates = @{
    "D2" = "246.17"
    "D3" = "22.17"
    "D4" = "5.356"
    "D5" = "0.05856"
    "D6" = "3.394"
    "D7" = "6.376"
    "D9" = "1.013"
    "D10" = "0.1422"
    "D11" = "0.04155"
    "D12" = "0.07349"
    "D13" = "0.03009"
    "D14" = "4.178"
    "D15" = "0.09408"
    "D16" = "0.001595"
    "D17" = "0.04810"
    "D19" = "0.005951"
    "D20" = "0.004084"
    "D21" = "0.0009871"
    "D23" = "3.689"
    "D41" = "0.006381"
    "D42" = "0.1073"
    "D43" = "0.003000"
    "D44" = "0.005127"
    "D45" = "0.00005656"
    "D48" = "0.08788"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.ClearFormats()
}
